# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "541.53") are assigned with a
# leading apostrophe so Excel stores them as literal text instead of converting them
# to a floating point number (which would corrupt values such as "8.50" -> 8.5).

$ws.Range("D2").Value = '58.210.81'
$ws.Range("E2").Value = '  +1.89%  '
$ws.Range("D3").Value = '2.358.72'
$ws.Range("E3").Value = '  +1.90%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''541.53'
$ws.Range("E5").Value = '  +1.71%  '
$ws.Range("D6").Value = '''135.44'
$ws.Range("E6").Value = '  +2.41%  '
$ws.Range("E7").Value = '  +0.61%  '
$ws.Range("D8").Value = '''0.564'
$ws.Range("E8").Value = '  +5.06%  '
$ws.Range("E9").Value = '  +0.71%  '
$ws.Range("D10").Value = '''5.61'
$ws.Range("E10").Value = '  +5.84%  '
$ws.Range("E11").Value = '  -0.57%  '
$ws.Range("D12").Value = '''0.355'
$ws.Range("E12").Value = '  +2.68%  '
$ws.Range("D13").Value = '''23.89'
$ws.Range("E13").Value = '  +2.00%  '
$ws.Range("D14").Value = '2.777.86'
$ws.Range("E14").Value = '  +1.62%  '
$ws.Range("D15").Value = '58.184.22'
$ws.Range("E15").Value = '  +1.78%  '
$ws.Range("D16").Value = '''0.0000133'
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("D17").Value = '2.343.46'
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("E18").Value = '  +3.19%  '
$ws.Range("D19").Value = '''332.60'
$ws.Range("E19").Value = '  -1.49%  '
$ws.Range("E20").Value = '  +3.08%  '
$ws.Range("D21").Value = '''6.85'
$ws.Range("E21").Value = '  -0.76%  '
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").Value = '''62.98'
$ws.Range("E23").Value = '  +2.14%  '
$ws.Range("D24").Value = '''0.167'
$ws.Range("E24").Value = '  +0.61%  '
$ws.Range("D25").Value = '''8.50'
$ws.Range("E25").Value = '  -2.22%  '
$ws.Range("E26").Value = '  +1.06%  '
$ws.Range("E27").Value = '  +5.38%  '
$ws.Range("E28").Value = '  +1.68%  '
$ws.Range("D29").Value = '''170.72'
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("E30").Value = '  +2.00%  '
$ws.Range("E31").Value = '  +1.09%  '
$ws.Range("D32").Value = '''18.47'
$ws.Range("E32").Value = '  -0.36%  '
$ws.Range("E33").Value = '  +12.69%  '
$ws.Range("D35").Value = '''4.27'
$ws.Range("E35").Value = '  +7.16%  '
$ws.Range("E36").Value = '  +0.78%  '
$ws.Range("E37").Value = '  -0.66%  '
$ws.Range("D38").Value = '''1.65'
$ws.Range("E38").Value = '  +3.60%  '
$ws.Range("D39").Value = '''39.17'
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").Value = '''147.37'
$ws.Range("E40").Value = '  -0.80%  '
$ws.Range("D41").Value = '''293.85'
$ws.Range("E41").Value = '  +4.50%  '
$ws.Range("D42").Value = '''0.377'
$ws.Range("E42").Value = '  +0.22%  '
$ws.Range("E43").Value = '  +1.42%  '
$ws.Range("E44").Value = '  +2.52%  '
$ws.Range("D45").Value = '''19.27'
$ws.Range("E45").Value = '  +3.08%  '
$ws.Range("D46").Value = '''0.0503'
$ws.Range("E46").Value = '  +0.47%  '
$ws.Range("D47").Value = '''0.563'
$ws.Range("E47").Value = '  +1.27%  '
$ws.Range("E48").Value = '  +1.51%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''17.53'
$ws.Range("E49").Value = '  +1.04%  '
$ws.Range("B50").Value = 'Polygon'
$ws.Range("C50").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D50").Value = '''0.382'
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("D51").Value = '''11.08'
$ws.Range("E51").Value = '  +0.62%  '
